# Insert a new data row at row 352 (pushing existing rows 352..460 down to
# 353..461) on the single worksheet of the workbook, then populate the new
# row with the record's values, matching the target OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 352:460 down to 353:461, creating a blank row 352.
$ws.Rows("352:352").Insert()

# Populate the newly inserted row 352 with the new record.
$ws.Range("A352").Value = 9
$ws.Range("B352").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C352").Value = "Metropolitana"
$ws.Range("D352").Value = 45229
$ws.Range("D352").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E352").Value = 13
$ws.Range("F352").Value = 100112043
$ws.Range("G352").Value = "Pepino ensalada"
$ws.Range("H352").Value = "Sin especificar"
$ws.Range("I352").Value = "Primera"
$ws.Range("J352").Value = 70
$ws.Range("K352").Value = 14000
$ws.Range("L352").Value = 15000
$ws.Range("M352").Value = 14500
$ws.Range("N352").Value = "`$/caja 60 unidades"
$ws.Range("O352").Value = "Región de Arica y Parinacota"
$ws.Range("P352").Value = 242
$ws.Range("Q352").Value = 60
$ws.Range("R352").Value = "Hortaliza"
